$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C115").Value = 56.6
$ws.Range("C183").Value = 5
$ws.Range("C185").Value = 60.3
$ws.Range("C186").Value = 50.8
$ws.Range("C187").Value = 61.4
$ws.Range("C188").Value = 92
$ws.Range("C190").Value = 8.1
$ws.Range("C191").Value = 92.8
$ws.Range("C266").Value = 96.59999999999999
$ws.Range("C267").Value = 46.7
$ws.Range("C268").Value = 35.2
$ws.Range("C269").Value = 30.9
$ws.Range("C271").Value = 97.3
$ws.Range("C272").Value = 62.7
$ws.Range("C273").Value = 99
$ws.Range("C274").Value = 44.6
$ws.Range("C275").Value = 63.5
$ws.Range("C277").Value = 63.7
$ws.Range("C365").Value = 33.6
$ws.Range("C424").Value = 74.2
$ws.Range("C425").Value = 97.90000000000001
$ws.Range("C426").Value = 88.40000000000001
$ws.Range("C431").Value = 130
